$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43-122 down to 44-123
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new data record
$ws.Cells.Item(43, 1).Value = 10
$ws.Cells.Item(43, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43, 3).Value = "La Araucanía"
$ws.Cells.Item(43, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(43, 5).Value = 9
$ws.Cells.Item(43, 6).Value = 300000000
$ws.Cells.Item(43, 7).Value = "Espárragos"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 30
$ws.Cells.Item(43, 11).Value = 1800
$ws.Cells.Item(43, 12).Value = 1800
$ws.Cells.Item(43, 13).Value = 1800
$ws.Cells.Item(43, 14).Value = "$/kilo"
$ws.Cells.Item(43, 15).Value = "Región del Maule"
$ws.Cells.Item(43, 16).Value = 1800
$ws.Cells.Item(43, 17).Value = 1
$ws.Cells.Item(43, 18).Value = "Hortaliza"
